# Update the "Time:" stamp inside the statsmodels OLS Regression Results
# summary text block that each backward-elimination worksheet carries in
# cell B2. All but the very last worksheet ("8") move to 20:59:48; the
# last one (the smallest / final model) moves to 20:59:49, matching the
# author's re-run before starting the documentation pass.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -like "*Time:*") {
        if ($i -eq $sheetCount) {
            $newTime = "20:59:49"
        } else {
            $newTime = "20:59:48"
        }

        $newText = $text -replace "(Time:\s+)\d\d:\d\d:\d\d", "`${1}$newTime"
        $cell.Value = $newText
    }
}
